$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.05261753845416351
$ws.Range("C2").Value = 1.350476871904311
$ws.Range("D2").Value = 10.10003955787528
$ws.Range("E2").Value = 3.178055940016677
$ws.Range("F2").Value = 3.209239220060225
$ws.Range("G2").Value = 51

# Row 3
$ws.Range("B3").Value = 0.02310908688600665
$ws.Range("C3").Value = 1.508687381939164
$ws.Range("D3").Value = 9.675801997808875
$ws.Range("E3").Value = 3.110595119556526
$ws.Range("F3").Value = 3.1420888619256
$ws.Range("G3").Value = 50

# Row 4
$ws.Range("B4").Value = 0.0619550970533476
$ws.Range("C4").Value = 1.373983353632171
$ws.Range("D4").Value = 6.861338155303133
$ws.Range("E4").Value = 2.619415613319722
$ws.Range("F4").Value = 2.645820156166766
$ws.Range("G4").Value = 49

# Row 5
$ws.Range("B5").Value = 0.06730251648952341
$ws.Range("C5").Value = 1.433898174105932
$ws.Range("D5").Value = 8.49624289535795
$ws.Range("E5").Value = 2.914831538075219
$ws.Range("F5").Value = 2.944891851492343
$ws.Range("G5").Value = 48

# Row 6
$ws.Range("B6").Value = 0.04899455259774384
$ws.Range("C6").Value = 1.583601773355311
$ws.Range("D6").Value = 9.97197792596331
$ws.Range("E6").Value = 3.157843872955614
$ws.Range("F6").Value = 3.19159950566434
$ws.Range("G6").Value = 47

# Row 7
$ws.Range("B7").Value = 0.05601300463246239
$ws.Range("C7").Value = 1.591143510471736
$ws.Range("D7").Value = 8.723097247822707
$ws.Range("E7").Value = 2.953488995717219
$ws.Range("F7").Value = 2.985588162386466
$ws.Range("G7").Value = 46

# Row 8
$ws.Range("B8").Value = 0.01539162827235637
$ws.Range("C8").Value = 1.578671680484494
$ws.Range("D8").Value = 8.737106999296079
$ws.Range("E8").Value = 2.955859773280201
$ws.Range("F8").Value = 2.989219852495662
$ws.Range("G8").Value = 45
